$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4-6 (the extra duplicate sample rows), working bottom-up
# so the remaining row indices don't shift out from under us.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Insert a new column B ("model") before the existing "prompt" column,
# shifting prompt/user_input/assistant_response/response_time right by one.
$ws.Columns.Item(2).Insert()

# Header row
$ws.Range("B1").Value = "model"

# Data rows
$ws.Range("B2").Value = "openai:gpt-4o-mini"
$ws.Range("B3").Value = "openai:gpt-4o-mini"

# Updated response_time values (now in column F after the insert)
$ws.Range("F2").Value = 5.405269861221313
$ws.Range("F3").Value = 0.5770392417907715
